$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02
    "C2" = 1.037494470852308
    "D2" = 1.046758875721007
    "E2" = 1.036277335339771
    "F2" = 1.055552244180352
    "I2" = 1.042222328563286
    "J2" = 1.0425970358803
    "K2" = 1.049523413466975
    "L2" = 1.039071536915578
    "M2" = 1.058292388311089
    "N2" = 1.044077643258801
    "B3" = 1.019999999999999
    "C3" = 1.038402555501242
    "D3" = 1.047520277088291
    "E3" = 1.037048448028472
    "F3" = 1.056545920864358
    "I3" = 1.042491551297867
    "J3" = 1.043149731014273
    "K3" = 1.050096652754297
    "L3" = 1.039652319972062
    "M3" = 1.059099069303024
    "N3" = 1.04463112328326
    "B4" = 1.02
    "C4" = 1.038990546374506
    "D4" = 1.048013261285008
    "E4" = 1.037548121978535
    "F4" = 1.057189728289448
    "I4" = 1.042664607748462
    "J4" = 1.0435071245227
    "K4" = 1.050467220535403
    "L4" = 1.040028177928662
    "M4" = 1.059621261534691
    "N4" = 1.044989024331462
    "B5" = 1.02
    "C5" = 1.039237832267341
    "D5" = 1.048220583886963
    "E5" = 1.037758354068949
    "F5" = 1.057460582842157
    "I5" = 1.042737085070832
    "J5" = 1.043657315004896
    "K5" = 1.050622920889353
    "L5" = 1.04018620012095
    "M5" = 1.059840841598431
    "N5" = 1.045139428101371
    "B6" = 1.02
    "C6" = 1.039279358168711
    "D6" = 1.048255398449506
    "E6" = 1.037793662840818
    "F6" = 1.057506072080936
    "I6" = 1.042749238152383
    "J6" = 1.043682529233357
    "K6" = 1.050649058571008
    "L6" = 1.040212733386432
    "M6" = 1.059877712981264
    "N6" = 1.045164678136928
    "B7" = 1.02
    "C7" = 1.03899385025062
    "D7" = 1.048016031258371
    "E7" = 1.037550930447093
    "F7" = 1.05719334668522
    "I7" = 1.042665577277716
    "J7" = 1.043509131603391
    "K7" = 1.050469301351699
    "L7" = 1.040030289385039
    "M7" = 1.059624195377501
    "N7" = 1.044991034262438
    "B8" = 1.02
    "C8" = 1.03780127896366
    "D8" = 1.047016130580299
    "E8" = 1.036537788114472
    "F8" = 1.055887888773962
    "I8" = 1.042313551113007
    "J8" = 1.042783870323763
    "K8" = 1.049717215539731
    "L8" = 1.039267803766078
    "M8" = 1.058564964555397
    "N8" = 1.044264743028604
    "B9" = 1.02
    "C9" = 1.035702922634822
    "D9" = 1.045256585788239
    "E9" = 1.034758020700707
    "F9" = 1.053593936418218
    "I9" = 1.041684465945159
    "J9" = 1.041504093079754
    "K9" = 1.048389260151612
    "L9" = 1.037924655304531
    "M9" = 1.056700168824138
    "N9" = 1.042983148354117
    "B10" = 1.02
    "C10" = 1.034306170526931
    "D10" = 1.044085254813632
    "E10" = 1.033575302830212
    "F10" = 1.052069035882431
    "I10" = 1.041259215647181
    "J10" = 1.040649772000852
    "K10" = 1.047502212285586
    "L10" = 1.037029586989564
    "M10" = 1.055458187297249
    "N10" = 1.042127614041291
    "B11" = 1.02
    "C11" = 1.033701884887305
    "D11" = 1.043578475158091
    "E11" = 1.033064090142856
    "F11" = 1.051409796887199
    "I11" = 1.041073695109829
    "J11" = 1.040279583769193
    "K11" = 1.047117709553464
    "L11" = 1.036642111837443
    "M11" = 1.054920698329651
    "N11" = 1.041756900099882
    "B12" = 1.02
    "C12" = 1.033477505003721
    "D12" = 1.043390298200932
    "E12" = 1.032874341386426
    "F12" = 1.051165085617068
    "I12" = 1.041004577107037
    "J12" = 1.040142040993942
    "K12" = 1.046974828222714
    "L12" = 1.036498201658776
    "M12" = 1.054721096774097
    "N12" = 1.041619161998113
    "B13" = 1.02
    "C13" = 1.033525631672601
    "D13" = 1.043430659899187
    "E13" = 1.032915036853773
    "F13" = 1.051217569779964
    "I13" = 1.041019412530077
    "J13" = 1.040171546108375
    "K13" = 1.047005479447962
    "L13" = 1.036529070169799
    "M13" = 1.054763909895684
    "N13" = 1.041648709013193
    "B14" = 1.02
    "C14" = 1.033683335964701
    "D14" = 1.043562919087673
    "E14" = 1.033048402620176
    "F14" = 1.05138956571743
    "I14" = 1.041067986022062
    "J14" = 1.040268215218202
    "K14" = 1.047105900154318
    "L14" = 1.036630215858231
    "M14" = 1.054904198262424
    "N14" = 1.041745515404244
    "B15" = 1.02
    "C15" = 1.033780513243609
    "D15" = 1.043644416801934
    "E15" = 1.03313059205094
    "F15" = 1.051495559256401
    "I15" = 1.041097886286153
    "J15" = 1.040327771188444
    "K15" = 1.047167764758142
    "L15" = 1.036692537119069
    "M15" = 1.054990640655043
    "N15" = 1.041805155950795
    "B16" = 1.02
    "C16" = 1.034346285957526
    "D16" = 1.04411889695317
    "E16" = 1.033609249664844
    "F16" = 1.052112809695006
    "I16" = 1.041271498916002
    "J16" = 1.040674334741686
    "K16" = 1.047527722012559
    "L16" = 1.037055304535131
    "M16" = 1.055493865027965
    "N16" = 1.042152211664034
    "B17" = 1.02
    "C17" = 1.034701319417554
    "D17" = 1.044416637687329
    "E17" = 1.033909744075416
    "F17" = 1.052500277382858
    "I17" = 1.041380031302106
    "J17" = 1.040891655518639
    "K17" = 1.047753405904379
    "L17" = 1.037282885196406
    "M17" = 1.05580960475117
    "N17" = 1.042369841061417
    "B18" = 1.02
    "C18" = 1.034908454267054
    "D18" = 1.044590344714295
    "E18" = 1.034085105278674
    "F18" = 1.052726382245362
    "I18" = 1.041443202798562
    "J18" = 1.041018389752292
    "K18" = 1.047885004343514
    "L18" = 1.037415638275628
    "M18" = 1.055993798992687
    "N18" = 1.042496755272219
    "B19" = 1.02
    "C19" = 1.034979090321098
    "D19" = 1.04464958104824
    "E19" = 1.034144913784001
    "F19" = 1.052803495341405
    "I19" = 1.04146471996684
    "J19" = 1.04106159853196
    "K19" = 1.047929869339197
    "L19" = 1.037460905162487
    "M19" = 1.05605660928644
    "N19" = 1.04254002541331
    "B20" = 1.02
    "C20" = 1.034663222524968
    "D20" = 1.044384688788395
    "E20" = 1.03387749476733
    "F20" = 1.052458695257657
    "I20" = 1.041368400619816
    "J20" = 1.040868341666124
    "K20" = 1.047729196201746
    "L20" = 1.037258466993509
    "M20" = 1.055775725900331
    "N20" = 1.042346494100557
    "B21" = 1.02
    "C21" = 1.03363689380056
    "D21" = 1.043523970287704
    "E21" = 1.033009125894449
    "F21" = 1.051338912789256
    "I21" = 1.041053688060427
    "J21" = 1.040239749626995
    "K21" = 1.047076330404817
    "L21" = 1.036600430542166
    "M21" = 1.054862885562265
    "N21" = 1.041717009388633
    "B22" = 1.02
    "C22" = 1.032992056079233
    "D22" = 1.042983170858987
    "E22" = 1.032463949329052
    "F22" = 1.050635784226753
    "I22" = 1.040854616027657
    "J22" = 1.039844307242213
    "K22" = 1.04666550170889
    "L22" = 1.036186786196464
    "M22" = 1.054289212030555
    "N22" = 1.041321005430305
    "B23" = 1.02
    "C23" = 1.033333853208861
    "D23" = 1.043269823504273
    "E23" = 1.032752881314253
    "F23" = 1.051008438064148
    "I23" = 1.040960261411427
    "J23" = 1.040053959433985
    "K23" = 1.046883322250554
    "L23" = 1.036406058085588
    "M23" = 1.054593301689559
    "N23" = 1.041530955352238
    "B24" = 1.02
    "C24" = 1.034680436716942
    "D24" = 1.044399124997465
    "E24" = 1.03389206657006
    "F24" = 1.052477484116592
    "I24" = 1.041373656437377
    "J24" = 1.040878876270792
    "K24" = 1.047740135645345
    "L24" = 1.037269500500219
    "M24" = 1.055791034206442
    "N24" = 1.042357043665573
    "B25" = 1.02
    "C25" = 1.036245023372289
    "D25" = 1.045711176702579
    "E25" = 1.035217470653604
    "F25" = 1.05418620802282
    "I25" = 1.041848135104598
    "J25" = 1.04183515059427
    "K25" = 1.048732880313363
    "L25" = 1.038271831787438
    "M25" = 1.057182054070627
    "N25" = 1.043314676008276
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Updated $($updates.Count) cells"